# "Generate Report for Archive"
# The handoff status text moves from "Ready for handoff" to "In Translation"
# on the Overview sheet (columns E/F, row 2) and on each per-language sheet
# (zh-cn, de-de; column C, row 2). Because the Status column is driven by
# this text, its column narrows to fit the new (shorter) label on every
# sheet that shows it.

$wb = $excel.ActiveWorkbook

# --- Update the status text everywhere it appears ---------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"

# --- Narrow the Status columns to match the shorter text ---------------
# (Column width in COM "characters" -- 12.5 is the value that resolves to
# the narrower rendered width used for the "In Translation" status text.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZh.Columns.Item(3).ColumnWidth = 12.5

$wsDe.Columns.Item(3).ColumnWidth = 12.5
